$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row strings: "_old" -> "_FV2410", "_new" -> "_FV2504"
$newHeadersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2410[$i]
}

$newHeadersFV2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2504[$i]
}

# 2) Convert the used range into an Excel Table ("Table1") with header row, default style
$tableRange = $ws.Range("A1:U64")
$listObj = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObj.Name = "Table1"

# 3) Freeze the header row (freeze panes at row 2 / split after row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
